$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Standard errors (theta_se in row 4, lambda_se in row 6) are now computed
# from the pickled bootstrap replicates instead of being left as "(nan)".
# Values are written column-by-column (theta_se then lambda_se for each
# year) so that the shared-string table is populated in the same order
# the values are produced by the replication pipeline.
$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L")

$row4 = @("(0.32)", "(1.62)", "(3.72)", "(2.81)", "(1.45)", "(3.46)", "(2.79)", "(2.56)", "(0.43)", "(3.1)", "(0.23)")
$row6 = @("(0.19)", "(0.47)", "(2.03)", "(1.04)", "(0.58)", "(2.85)", "(2.39)", "(1.63)", "(0.52)", "(1.95)", "(0.02)")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $col = $cols[$i]
    $ws.Range("$col" + "4").Value = $row4[$i]
    $ws.Range("$col" + "6").Value = $row6[$i]
}
